# Applies the "handles float input without breaking stuff" marksheet edit:
#  - recompute Right/Wrong/Not-Attempt/Max + Total row
#  - mark the "Marking" wrong-answer penalty as a real number instead of text
#  - collapse the extra "2nd"/"3rd" question-set columns (D/E beyond row 18,
#    and G/H entirely) which are no longer used
#  - fill in the "Student Ans" column for the questions the student got
#    right, with the "correctStyle" highlight

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- summary block (rows 10-12) ------------------------------------------

$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 7
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 21
$ws.Range("E10").Value = 28

$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 28
$ws.Range("C12").Value = 0
$ws.Range("E12").Value = "28/112"

# --- drop the unused G/H "3rd" answer columns -----------------------------

$unusedGH = @("G15","H15","G16","H16","G17","H17","G18","H18","G19","H19","G20","H20","G21","H21")
foreach ($addr in $unusedGH) {
    $ws.Range($addr).Clear()
}

# --- drop the unused D/E "2nd" answer column for rows 19-40 --------------

for ($r = 19; $r -le 40; $r++) {
    $ws.Range("D$r").Clear()
    $ws.Range("E$r").Clear()
}

# --- fill in "Student Ans" for the questions answered correctly ----------

$correct = @{
    20 = "Option B"
    22 = "Option D"
    23 = "Option D"
    27 = "Option A"
    35 = "Option D"
    37 = "Option A"
    38 = "Option A"
}
foreach ($r in $correct.Keys) {
    $ws.Range("A$r").Value = $correct[$r]
    $ws.Range("A$r").Style = "correctStyle"
}
